$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / title text updates ---
$ws.Range("M6").Value = "Thomas G. Donlon"
$ws.Range("A8").Value = "Volume 31   Number  39"
$ws.Range("C9").Value = "Report Covering the Week  9/23/2024  Through  9/29/2024"

# --- Cells changing between numeric and "0"-text representation ---
# C15, C29, C30: numeric count -> text "0" (reuse formatting+value from C20, an existing text-"0" cell)
$ws.Range("C20").Copy($ws.Range("C15"))
$ws.Range("C20").Copy($ws.Range("C29"))
$ws.Range("C20").Copy($ws.Range("C30"))
# F31: text "0" -> numeric 2 (reuse formatting+value from F29, an existing numeric "2" cell)
$ws.Range("F29").Copy($ws.Range("F31"))

# --- Weekly crime statistics data refresh ---
# Row 15
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = -100
$ws.Range("F15").Value = 4
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = 33.333333333333
$ws.Range("I15").Value = 23
$ws.Range("J15").Value = 13
$ws.Range("K15").Value = 76.923076923076
$ws.Range("L15").Value = 27.777777777777
$ws.Range("M15").Value = 155.555555555556
$ws.Range("N15").Value = -30.303030303030

# Row 16
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -60
$ws.Range("F16").Value = 28
$ws.Range("G16").Value = 27
$ws.Range("H16").Value = 3.703703703703
$ws.Range("I16").Value = 304
$ws.Range("J16").Value = 351
$ws.Range("K16").Value = -13.390313390313
$ws.Range("L16").Value = -34.199134199134
$ws.Range("M16").Value = 162.068965517241
$ws.Range("N16").Value = -84.313725490196

# Row 17
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 35
$ws.Range("G17").Value = 36
$ws.Range("H17").Value = -2.777777777777
$ws.Range("I17").Value = 402
$ws.Range("J17").Value = 378
$ws.Range("K17").Value = 6.349206349206
$ws.Range("L17").Value = 16.860465116279
$ws.Range("M17").Value = 179.166666666667
$ws.Range("N17").Value = -20.710059171597

# Row 18
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 60
$ws.Range("F18").Value = 21
$ws.Range("G18").Value = 34
$ws.Range("H18").Value = -38.235294117647
$ws.Range("I18").Value = 273
$ws.Range("J18").Value = 312
$ws.Range("K18").Value = -12.5
$ws.Range("L18").Value = -46.259842519685
$ws.Range("M18").Value = 9.638554216867
$ws.Range("N18").Value = -86.411149825784

# Row 19
$ws.Range("C19").Value = 47
$ws.Range("D19").Value = 43
$ws.Range("E19").Value = 9.302325581395
$ws.Range("F19").Value = 147
$ws.Range("G19").Value = 170
$ws.Range("H19").Value = -13.529411764705
$ws.Range("I19").Value = 1506
$ws.Range("J19").Value = 1745
$ws.Range("K19").Value = -13.696275071633
$ws.Range("L19").Value = -10.357142857142
$ws.Range("M19").Value = -11.097992916174
$ws.Range("N19").Value = -78.857223080162

# Row 20
$ws.Range("D20").Value = 1
$ws.Range("F20").Value = 2
$ws.Range("H20").Value = -80
$ws.Range("J20").Value = 54
$ws.Range("K20").Value = -22.222222222222
$ws.Range("L20").Value = -14.285714285714
$ws.Range("M20").Value = 133.333333333333
$ws.Range("N20").Value = -84.210526315789

# Row 21
$ws.Range("C21").Value = 63
$ws.Range("D21").Value = 65
$ws.Range("E21").Value = -3.076923076923
$ws.Range("F21").Value = 237
$ws.Range("G21").Value = 280
$ws.Range("H21").Value = -15.357142857142
$ws.Range("I21").Value = 2554
$ws.Range("J21").Value = 2855
$ws.Range("K21").Value = -10.542907180385
$ws.Range("L21").Value = -16.726442777958
$ws.Range("M21").Value = 14.529147982062
$ws.Range("N21").Value = -78.510727808161

# Row 22
$ws.Range("C22").Value = 3
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = 50
$ws.Range("F22").Value = 7
$ws.Range("H22").Value = -41.666666666666
$ws.Range("I22").Value = 121
$ws.Range("J22").Value = 157
$ws.Range("K22").Value = -22.929936305732
$ws.Range("L22").Value = -15.972222222222
$ws.Range("M22").Value = 12.037037037037

# Row 24
$ws.Range("C24").Value = 105
$ws.Range("D24").Value = 78
$ws.Range("E24").Value = 34.615384615384
$ws.Range("F24").Value = 349
$ws.Range("G24").Value = 304
$ws.Range("H24").Value = 14.802631578947
$ws.Range("I24").Value = 3285
$ws.Range("J24").Value = 3037
$ws.Range("K24").Value = 8.165953243332
$ws.Range("L24").Value = 32.888349514563
$ws.Range("M24").Value = -6.25

# Row 25
$ws.Range("C25").Value = 85
$ws.Range("D25").Value = 62
$ws.Range("E25").Value = 37.096774193548
$ws.Range("F25").Value = 299
$ws.Range("G25").Value = 246
$ws.Range("H25").Value = 21.544715447154
$ws.Range("I25").Value = 2891
$ws.Range("J25").Value = 2747
$ws.Range("K25").Value = 5.242082271568
$ws.Range("L25").Value = 25.805047867711

# Row 26
$ws.Range("C26").Value = 21
$ws.Range("E26").Value = 50
$ws.Range("F26").Value = 86
$ws.Range("G26").Value = 78
$ws.Range("H26").Value = 10.256410256410
$ws.Range("I26").Value = 785
$ws.Range("J26").Value = 802
$ws.Range("K26").Value = -2.119700748129
$ws.Range("L26").Value = 18.759455370650
$ws.Range("M26").Value = 79.223744292237

# Row 27
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 66.666666666666
$ws.Range("J27").Value = 22
$ws.Range("K27").Value = 31.818181818181

# Row 28
$ws.Range("C28").Value = 8
$ws.Range("D28").Value = 5
$ws.Range("E28").Value = 60
$ws.Range("F28").Value = 22
$ws.Range("H28").Value = 83.333333333333
$ws.Range("I28").Value = 176
$ws.Range("J28").Value = 166
$ws.Range("K28").Value = 6.024096385542
$ws.Range("L28").Value = 4.142011834319

# Row 31
$ws.Range("I31").Value = 19
$ws.Range("J31").Value = 12
$ws.Range("K31").Value = 58.333333333333
$ws.Range("L31").Value = -5
